# Atualiza os resultados da tabela "finalidade_meds":
#   n (Tratamento):        38        -> 37
#   SPT (%) (Tratamento):  4 (10.5)  -> 4 (10.8)
#   TEP (%) (Tratamento):  2 ( 5.3)  -> 2 ( 5.4)
#   TEP (%) (p):           0.361     -> 0.349

$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(2, 3).Range.Text = "37"
$t.Cell(3, 3).Range.Text = "4 (10.8)"
$t.Cell(4, 3).Range.Text = "2 ( 5.4)"
$t.Cell(4, 4).Range.Text = "0.349"
